# Apply the "Add Category related Menus" update to the GUI sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GUI")

# Mark the "ViewPersonalInfoMenu" task (row 29, assigned to Aryan) as done.
$ws.Range("C29").Value = 1

# Recalculate the workbook so dependent formulas (G14, H14) update.
$excel.Calculate()

# Move the active selection to C1, as it was left after the edit.
$ws.Activate()
$ws.Range("C1").Select()
